$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append two new batches of 5 device-master rows (rows 147-156), mirroring
# the existing repeating pattern of:
#   Finger Print Scanner / IRIS Scanner / Web Camera / Document Scanner / Printer
# ---------------------------------------------------------------------------

# Batch 1 -> rows 147-151
$batch1Ids    = @(3000166, 3000167, 3000168, 3000169, 3000170)
$batch1Names  = @("Finger Print Scanner 30", "IRIS Scanner 30", "Web Camera 30", "Document Scanner 30", "Printer 30")
$batch1Macs   = @("D6-15-AC-80-6B-86", "6D-58-E2-DF-74-34", "E2-A8-56-86-15-30", "72-E8-B9-FD-63-65", "D3-F3-A4-50-AD-12")
$batch1Serial = @("BS563Q2230814", "BS563Q2230815", "BS563Q2230816", "BS563Q2230817", "BS563Q2230818")
$batch1Dspec  = @(165, 327, 736, 801, 920)

# Batch 2 -> rows 152-156
$batch2Ids    = @(3000171, 3000172, 3000173, 3000174, 3000175)
$batch2Names  = @("Finger Print Scanner 31", "IRIS Scanner 31", "Web Camera 31", "Document Scanner 31", "Printer 31")
$batch2Macs   = @("06-16-D0-0B-A6-E4", "21-78-45-AC-E9-20", "3C-E8-87-99-DB-FA", "BF-55-53-98-40-08", "5A-43-36-46-22-EB")
$batch2Serial = @("BS563Q2230819", "BS563Q2230820", "BS563Q2230821", "BS563Q2230822", "BS563Q2230823")
$batch2Dspec  = @(165, 327, 736, 801, 920)

$batch1Start = 147
$batch2Start = 152

# --- Batch 1: ids, then names (col B), then macs (col C), then serials (col D) ---
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch1Start + $i, 1).Value = $batch1Ids[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch1Start + $i, 2).Value = $batch1Names[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch1Start + $i, 3).Value = $batch1Macs[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch1Start + $i, 4).Value = $batch1Serial[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch1Start + $i, 6).Value = $batch1Dspec[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch1Start + $i, 7).Value = "eng"
}
for ($i = 0; $i -lt 5; $i++) {
    $hCell = $ws.Cells.Item($batch1Start + $i, 8)
    $hCell.Value = $true
    $hCell.HorizontalAlignment = -4131
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch1Start + $i, 9).Value = "superadmin"
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch1Start + $i, 10).Value = "now()"
}

# --- Batch 2: ids, then names (col B), then serials (col D), then macs (col C) ---
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch2Start + $i, 1).Value = $batch2Ids[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch2Start + $i, 2).Value = $batch2Names[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch2Start + $i, 4).Value = $batch2Serial[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch2Start + $i, 3).Value = $batch2Macs[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch2Start + $i, 6).Value = $batch2Dspec[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch2Start + $i, 7).Value = "eng"
}
for ($i = 0; $i -lt 5; $i++) {
    $hCell = $ws.Cells.Item($batch2Start + $i, 8)
    $hCell.Value = $true
    $hCell.HorizontalAlignment = -4131
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch2Start + $i, 9).Value = "superadmin"
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($batch2Start + $i, 10).Value = "now()"
}

# Update the visible selection to match the post-edit state.
$ws.Range("E156").Select() | Out-Null
